# Refresh the cryptos table (prices + 1h volume deltas) with the
# latest scrape, matching the GitHub Actions bot's daily commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") stores figures as text in the source data (values like
# "64.100.00" use "." as a thousands separator and are not valid numbers), so
# numeric-looking replacements are written with a leading apostrophe to force
# Excel to keep them as text instead of auto-converting to a number.

$ws.Range("D2").Value = '64.098.38'
$ws.Range("E2").Value = '  +0.19%  '

$ws.Range("D3").Value = '2.760.05'

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").Value = "'" + '579.27'
$ws.Range("E5").Value = '  +0.67%  '

$ws.Range("D6").Value = "'" + '158.75'
$ws.Range("E6").Value = '  +2.74%  '

$ws.Range("E7").Value = '  +0.18%  '

$ws.Range("D8").Value = "'" + '0.609'
$ws.Range("E8").Value = '  +0.19%  '

$ws.Range("E9").Value = '  -1.17%  '

$ws.Range("D10").Value = "'" + '5.74'
$ws.Range("E10").Value = '  -14.17%  '

$ws.Range("E11").Value = '  -0.62%  '

$ws.Range("D12").Value = "'" + '0.159'
$ws.Range("E12").Value = '  -2.27%  '

$ws.Range("D13").Value = '3.248.51'
$ws.Range("E13").Value = '  +1.18%  '

$ws.Range("D14").Value = "'" + '27.01'
$ws.Range("E14").Value = '  +2.69%  '

$ws.Range("D15").Value = '63.771.75'
$ws.Range("E15").Value = '  -0.06%  '

$ws.Range("E16").Value = '  +0.39%  '

$ws.Range("D17").Value = '2.761.98'
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("E18").Value = '  +1.97%  '

$ws.Range("E19").Value = '  +0.70%  '

$ws.Range("D20").Value = "'" + '361.06'
$ws.Range("E20").Value = '  +0.26%  '

$ws.Range("E21").Value = '  -1.14%  '

$ws.Range("D22").Value = "'" + '0.550'
$ws.Range("E22").Value = '  +3.17%  '

$ws.Range("E23").Value = '  +0.31%  '

$ws.Range("D24").Value = "'" + '65.83'
$ws.Range("E24").Value = '  -0.30%  '

$ws.Range("E25").Value = '  +1.51%  '

$ws.Range("D26").Value = "'" + '8.54'
$ws.Range("E26").Value = '  +0.32%  '

$ws.Range("E27").Value = '  +0.17%  '

$ws.Range("D28").Value = '0.0₃0930'
$ws.Range("E28").Value = '  +2.53%  '

$ws.Range("D29").Value = "'" + '1.97'
$ws.Range("E29").Value = '  -1.61%  '

$ws.Range("E30").Value = '  -0.87%  '

$ws.Range("E31").Value = '  +1.18%  '

$ws.Range("D32").Value = "'" + '167.56'
$ws.Range("E32").Value = '  -2.12%  '

$ws.Range("E33").Value = '  -0.54%  '

$ws.Range("E34").Value = '  +3.90%  '

$ws.Range("E35").Value = '  +0.20%  '

$ws.Range("E36").Value = '  +2.41%  '

$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("D38").Value = "'" + '0.994'
$ws.Range("E38").Value = '  -0.15%  '

$ws.Range("D39").Value = "'" + '6.27'
$ws.Range("E39").Value = '  +12.30%  '

$ws.Range("E40").Value = '  -0.88%  '

$ws.Range("D41").Value = "'" + '331.15'
$ws.Range("E41").Value = '  -4.00%  '

$ws.Range("D42").Value = "'" + '39.40'
$ws.Range("E42").Value = '  +0.22%  '

$ws.Range("D43").Value = "'" + '21.69'
$ws.Range("E43").Value = '  -0.14%  '

$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").Value = "'" + '0.0597'
$ws.Range("E44").Value = '  +1.17%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = "'" + '21.85'
$ws.Range("E45").Value = '  +0.28%  '

$ws.Range("D46").Value = "'" + '0.0258'
$ws.Range("E46").Value = '  +1.17%  '

$ws.Range("D48").Value = "'" + '136.76'
$ws.Range("E48").Value = '  -1.69%  '

$ws.Range("E49").Value = '  +0.90%  '

$ws.Range("E50").Value = '  +0.11%  '

$ws.Range("D51").Value = "'" + '11.06'
